$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 35, shifting existing rows 35..120 down to 36..121
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new data record
$ws.Cells.Item(35, 1).Value = 10
$ws.Cells.Item(35, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(35, 3).Value = "La Araucanía"
$ws.Cells.Item(35, 4).Value = 44914
$ws.Cells.Item(35, 5).Value = 9
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100101
$ws.Cells.Item(35, 8).Value = "Berries"
$ws.Cells.Item(35, 9).Value = 100101001
$ws.Cells.Item(35, 10).Value = "Arándano (blue)"
$ws.Cells.Item(35, 11).Value = "Sin especificar"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 580
$ws.Cells.Item(35, 14).Value = 1800
$ws.Cells.Item(35, 15).Value = 2000
$ws.Cells.Item(35, 16).Value = 1869
$ws.Cells.Item(35, 17).Value = "`$/kilo"
$ws.Cells.Item(35, 18).Value = "Región del Maule"
$ws.Cells.Item(35, 19).Value = 1869
$ws.Cells.Item(35, 20).Value = 1
